$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Current Fringe Amount" (GA fringe rate) from 0.1% to 0.6%
$ws.Range("I11").Value = 0.006

# Update the "Current PI Fringe Amount" (PI fringe rate) from 17.06% to 18.54%
$ws.Range("I18").Value = 0.1854
$ws.Range("J18").Value = 0.1854

# Move the active selection to E11 (matches the author's cursor position when saving)
$ws.Range("E11").Select()
